$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Rbp4 -> Stra6 -> FAPs (Sending/Target same as before, but now part of new set)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rbp4"
$ws.Range("C2").Value = "Stra6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.653109333333333
$ws.Range("H2").Value = 13.959328
$ws.Range("I2").Value = 0.981686136247225
$ws.Range("J2").Value = 0.9867633043798142
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.339544
$ws.Range("N2").Value = 1.018632
$ws.Range("O2").Value = 0.1965227220570919
$ws.Range("P2").Value = 0.2684097734345912
$ws.Range("Q2").Value = 1.579935355477333
$ws.Range("R2").Value = 14.219418199296
$ws.Range("S2").Value = 0.1929236317010138
$ws.Range("T2").Value = 0.2648569149621545

# Row 3: FAPs -> Rbp4 -> Stra6 -> sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rbp4"
$ws.Range("C3").Value = "Stra6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.653109333333333
$ws.Range("H3").Value = 13.959328
$ws.Range("I3").Value = 0.981686136247225
$ws.Range("J3").Value = 0.9867633043798142
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.3882155
$ws.Range("N3").Value = 2.776431
$ws.Range("O3").Value = 0.8034772779429081
$ws.Range("P3").Value = 0.7315902265654087
$ws.Range("Q3").Value = 6.459518499727999
$ws.Range("R3").Value = 38.75711099836799
$ws.Range("S3").Value = 0.7887625045462112
$ws.Range("T3").Value = 0.7219063894176596

# Row 4: Neutro -> Rbp4 -> Stra6 -> FAPs
$ws.Range("A4").Value = "Neutro"
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01364166666666667
$ws.Range("H4").Value = 0.040925
$ws.Range("I4").Value = 0.002878040055074119
$ws.Range("J4").Value = 0.002892924948231312
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.339544
$ws.Range("N4").Value = 1.018632
$ws.Range("O4").Value = 0.1965227220570919
$ws.Range("P4").Value = 0.2684097734345912
$ws.Range("Q4").Value = 0.004631946066666667
$ws.Range("R4").Value = 0.0416875146
$ws.Range("S4").Value = 0.0005656002658125085
$ws.Range("T4").Value = 0.000776489329918043

# Row 5: Neutro -> Rbp4 -> Stra6 -> sCs
$ws.Range("A5").Value = "Neutro"
$ws.Range("B5").Value = "Rbp4"
$ws.Range("C5").Value = "Stra6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01364166666666667
$ws.Range("H5").Value = 0.040925
$ws.Range("I5").Value = 0.002878040055074119
$ws.Range("J5").Value = 0.002892924948231312
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.3882155
$ws.Range("N5").Value = 2.776431
$ws.Range("O5").Value = 0.8034772779429081
$ws.Range("P5").Value = 0.7315902265654087
$ws.Range("Q5").Value = 0.0189375731125
$ws.Range("R5").Value = 0.113625438675
$ws.Range("S5").Value = 0.00231243978926161
$ws.Range("T5").Value = 0.002116435618313269

# Row 6: sCs -> Rbp4 -> Stra6 -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rbp4"
$ws.Range("C6").Value = "Stra6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.07316449999999999
$ws.Range("H6").Value = 0.146329
$ws.Range("I6").Value = 0.01543582369770094
$ws.Range("J6").Value = 0.01034377067195454
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.339544
$ws.Range("N6").Value = 1.018632
$ws.Range("O6").Value = 0.1965227220570919
$ws.Range("P6").Value = 0.2684097734345912
$ws.Range("Q6").Value = 0.024842566988
$ws.Range("R6").Value = 0.149055401928
$ws.Range("S6").Value = 0.003033490090265554
$ws.Range("T6").Value = 0.002776369142518687

# Row 7: sCs -> Rbp4 -> Stra6 -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rbp4"
$ws.Range("C7").Value = "Stra6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.07316449999999999
$ws.Range("H7").Value = 0.146329
$ws.Range("I7").Value = 0.01543582369770094
$ws.Range("J7").Value = 0.01034377067195454
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.3882155
$ws.Range("N7").Value = 2.776431
$ws.Range("O7").Value = 0.8034772779429081
$ws.Range("P7").Value = 0.7315902265654087
$ws.Range("Q7").Value = 0.10156809294975
$ws.Range("R7").Value = 0.4062723717989999
$ws.Range("S7").Value = 0.01240233360743539
$ws.Range("T7").Value = 0.007567401529435852
